# Insert a new row above row 38 (shifts existing rows 38-106 down to 39-107)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with data.
# Most fields mirror the old row 38 (now row 39) except date, variety,
# volume, weighted price, origin and price/kg which change.
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value = "Arica y Parinacota"
$ws.Range("D38").Value = 44868
$ws.Range("D38").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100102
$ws.Range("H38").Value = "Cítricos"
$ws.Range("I38").Value = 100102005
$ws.Range("J38").Value = "Naranja"
$ws.Range("K38").Value = "Lane Late"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 350
$ws.Range("N38").Value = 650
$ws.Range("O38").Value = 700
$ws.Range("P38").Value = 686
$ws.Range("Q38").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R38").Value = "Región de Coquimbo"
$ws.Range("S38").Value = 686
$ws.Range("T38").Value = 1
